$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Val)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "65.031.81"
Set-TextValue $ws "E2" "  -2.06%  "
Set-TextValue $ws "D3" "3.412.98"
Set-TextValue $ws "E3" "  -5.20%  "
Set-TextValue $ws "D4" "0.998"
Set-TextValue $ws "E4" "  -0.28%  "
Set-TextValue $ws "D5" "591.59"
Set-TextValue $ws "E5" "  -2.55%  "
Set-TextValue $ws "D6" "134.38"
Set-TextValue $ws "E6" "  -9.49%  "
Set-TextValue $ws "D7" "3.413.51"
Set-TextValue $ws "E7" "  -5.13%  "
Set-TextValue $ws "D8" "0.998"
Set-TextValue $ws "E8" "  -0.26%  "
Set-TextValue $ws "D9" "0.487"
Set-TextValue $ws "E9" "  +0.11%  "
Set-TextValue $ws "D10" "7.40"
Set-TextValue $ws "E10" "  -6.16%  "
Set-TextValue $ws "D11" "0.119"
Set-TextValue $ws "E11" "  -12.20%  "
Set-TextValue $ws "D12" "0.374"
Set-TextValue $ws "E12" "  -9.92%  "
Set-TextValue $ws "D13" "3.973.70"
Set-TextValue $ws "E13" "  -5.66%  "
Set-TextValue $ws "D14" "0.0000177"
Set-TextValue $ws "E14" "  -13.85%  "
Set-TextValue $ws "B15" "WrappedEther"
Set-TextValue $ws "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D15" "3.406.54"
Set-TextValue $ws "E15" "  -5.27%  "
Set-TextValue $ws "B16" "Avalanche"
Set-TextValue $ws "C16" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws "D16" "26.17"
Set-TextValue $ws "E16" "  -11.70%  "
Set-TextValue $ws "B17" "TRON"
Set-TextValue $ws "C17" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D17" "0.114"
Set-TextValue $ws "E17" "  -2.88%  "
Set-TextValue $ws "D18" "64.850.48"
Set-TextValue $ws "E18" "  -2.47%  "
Set-TextValue $ws "D19" "10.02"
Set-TextValue $ws "E19" "  -10.55%  "
Set-TextValue $ws "D20" "5.68"
Set-TextValue $ws "E20" "  -10.19%  "
Set-TextValue $ws "D21" "13.57"
Set-TextValue $ws "E21" "  -9.13%  "
Set-TextValue $ws "D22" "389.46"
Set-TextValue $ws "E22" "  -8.55%  "
Set-TextValue $ws "B23" "Litecoin"
Set-TextValue $ws "C23" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws "D23" "72.67"
Set-TextValue $ws "E23" "  -7.45%  "
Set-TextValue $ws "B24" "Polygon"
Set-TextValue $ws "C24" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws "D24" "0.541"
Set-TextValue $ws "E24" "  -11.69%  "
Set-TextValue $ws "E25" "  +0.08%  "
Set-TextValue $ws "D26" "3.551.21"
Set-TextValue $ws "E26" "  -5.07%  "
Set-TextValue $ws "D27" "0.0000103"
Set-TextValue $ws "E27" "  -14.67%  "
Set-TextValue $ws "B28" "Binance-PegBSC-USD"
Set-TextValue $ws "C28" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws "D28" "0.999"
Set-TextValue $ws "E28" "  -0.05%  "
Set-TextValue $ws "B29" "PancakeSwap"
Set-TextValue $ws "C29" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D29" "2.23"
Set-TextValue $ws "E29" "  -10.80%  "
Set-TextValue $ws "D30" "7.06"
Set-TextValue $ws "E30" "  -14.99%  "
Set-TextValue $ws "D31" "8.06"
Set-TextValue $ws "E31" "  -14.03%  "
Set-TextValue $ws "D32" "3.413.40"
Set-TextValue $ws "E32" "  -5.13%  "
Set-TextValue $ws "E33" "  -0.02%  "
Set-TextValue $ws "D34" "0.140"
Set-TextValue $ws "E34" "  -10.98%  "
Set-TextValue $ws "D35" "22.30"
Set-TextValue $ws "E35" "  -11.37%  "
Set-TextValue $ws "B36" "Monero"
Set-TextValue $ws "C36" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D36" "173.08"
Set-TextValue $ws "E36" "  -0.82%  "
Set-TextValue $ws "B37" "Fetch.AI"
Set-TextValue $ws "C37" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws "D37" "1.22"
Set-TextValue $ws "E37" "  -15.17%  "
Set-TextValue $ws "D38" "6.75"
Set-TextValue $ws "E38" "  -13.03%  "
Set-TextValue $ws "D39" "1.51"
Set-TextValue $ws "E39" "  -9.67%  "
Set-TextValue $ws "D40" "4.73"
Set-TextValue $ws "E40" "  -15.74%  "
Set-TextValue $ws "D41" "0.0755"
Set-TextValue $ws "E41" "  -10.99%  "
Set-TextValue $ws "D42" "43.84"
Set-TextValue $ws "E42" "  -4.52%  "
Set-TextValue $ws "D43" "0.805"
Set-TextValue $ws "E43" "  -9.46%  "
Set-TextValue $ws "B44" "FirstDigitalUSD"
Set-TextValue $ws "C44" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D44" "0.997"
Set-TextValue $ws "E44" "  -0.31%  "
Set-TextValue $ws "B45" "Filecoin"
Set-TextValue $ws "C45" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D45" "4.37"
Set-TextValue $ws "E45" "  -15.88%  "
Set-TextValue $ws "D46" "1.59"
Set-TextValue $ws "E46" "  -14.15%  "
Set-TextValue $ws "D47" "1.05"
Set-TextValue $ws "E47" "  -7.28%  "
Set-TextValue $ws "B48" "EnergySwap"
Set-TextValue $ws "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D48" "21.39"
Set-TextValue $ws "E48" "  -9.51%  "
Set-TextValue $ws "B49" "Cosmos"
Set-TextValue $ws "C49" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D49" "6.49"
Set-TextValue $ws "E49" "  -9.22%  "
Set-TextValue $ws "D50" "2.10"
Set-TextValue $ws "E50" "  -17.80%  "
Set-TextValue $ws "D51" "2.186.83"
Set-TextValue $ws "E51" "  -8.79%  "

Write-Host "Applied cryptos update"
